$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row labels (persistent URI refactoring: snake_case headers)
$ws.Range("A1").Value = "codice_1_livello"
$ws.Range("B1").Value = "label _ITA _1 _livello"
$ws.Range("C1").Value = "label_ENG_1_livello"
$ws.Range("D1").Value = "definizione"

# Update the view: scroll so column B is the left-most visible column,
# and move the selection to D1
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D1").Select()
